$d = $word.ActiveDocument

# --- Part 1: extend the "En la técnica de personas..." paragraph ---
# Append the new trailing sentence right after "...elementos de ésta."
$r = $d.Content
$r.Find.Execute("elementos de ésta.") | Out-Null
$r.InsertAfter(" Dicha información se obtendrá de los datos de uso que Facebook proporciona de cada una de las personas que decidan usar la aplicación, en base a las investigaciones que se tienen y las que se sumarán en un futuro próximo, se irán armando las ")

# "personas" in italics
$pos = $r.End
$rItalic = $d.Range($pos, $pos)
$rItalic.InsertAfter("personas")
$rItalic.Italic = 1

# trailing closing text (non-italic) + language tag
$pos2 = $rItalic.End
$rTail = $d.Range($pos2, $pos2)
$rTail.InsertAfter(" para tener la interfaz que corresponda para cada una.")
$rTail.LanguageID = "es-ES"

# --- Part 2: append two new paragraphs after "Debido al contexto..." ---
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$emptyPara = $d.Paragraphs.Last
$emptyPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "Aunque de ser necesario podríamos recabar datos sobre algún estudio médico, en caso de existir, donde se haya comprobado que las personas en la situación con la que vamos a tratar, tengan alguna cierta orientación hacia algún tipo de patrones o interfaz diferente, pero esto se llevará conforme al paso de los datos recabados por las investigacionesX"

# move the _GoBack bookmark from "Debido al contexto..." paragraph to the end
# of the new last paragraph (collapsed, zero-length, right after the text).
$d.Bookmarks("_GoBack").Delete()
$bmPos = $d.Content.End - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# drop the temporary trailing placeholder character used to dodge the
# end-of-document zero-length-range quirk
$tailStart = $d.Content.End - 2
$tailEnd = $d.Content.End - 1
$tailRange = $d.Range($tailStart, $tailEnd)
$tailRange.Text = ""
